$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Inscricoes")

# Row 8: Inscritos 12 -> 13
$ws.Range("E8").Value = 13

# Row 16: Inscritos 9 -> 10
$ws.Range("E16").Value = 10

# Row 28: Inscritos 13 -> 14
$ws.Range("E28").Value = 14

# Row 60: Inscritos 16 -> 17
$ws.Range("E60").Value = 17

# Row 71: Inscritos 26 -> 27, Pagos 12 -> 13, Inscricoes homologadas 12 -> 13
$ws.Range("E71").Value = 27
$ws.Range("F71").Value = 13
$ws.Range("H71").Value = 13

# Row 77: Inscritos 48 -> 49
$ws.Range("E77").Value = 49

# Row 89: Inscritos 28 -> 29
$ws.Range("E89").Value = 29
